$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) for the two adversarial-training model columns ---
$ws.Range("D1").Value = "adv_training_model(eps=0.2,alpha = 0.5)"
$ws.Range("E1").Value = "adv_training_model(eps=0.1,alpha = 0.5)"

# --- New result values for rows 2-5 (clean, eps=0.1, eps=0.2, eps=0.3) ---
$ws.Range("D2").Value = 0.97
$ws.Range("E2").Value = 0.97

$ws.Range("D3").Value = 0.97
$ws.Range("E3").Value = 0.96

$ws.Range("D4").Value = 0.97
$ws.Range("E4").Value = 0.92

$ws.Range("D5").Value = 0.95
$ws.Range("E5").Value = 0.81

# --- Column widths for the new columns (match the bestFit-style wide columns) ---
$ws.Columns.Item(4).ColumnWidth = 37.5
$ws.Columns.Item(5).ColumnWidth = 37.5

# --- Selection moves to G4 after the edit ---
$ws.Range("G4").Select()
